$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing shared string (row 22, col I) to add the rr_add suffix
$ws.Cells.Item(22, 9).Value = 'I/O bound {\tt [$n$: 8] [seed: 64] [$\lambda$: 0.001] [limit: 4096] [$t_{cs}$: 20] [$\alpha$: 0.5] [$t_{slice}$: 2048] [rr_{add}: BEGINNING]}'

# Add new data rows (Run Number 6-10)
$ws.Cells.Item(27, 1).Value = 6
$ws.Cells.Item(27, 2).Value = "FCFS"
$ws.Cells.Item(27, 3).Value = 1629.22
$ws.Cells.Item(27, 4).Value = 3216.87
$ws.Cells.Item(27, 5).Value = 4850.09
$ws.Cells.Item(27, 6).Value = 23
$ws.Cells.Item(27, 7).Value = 0
$ws.Cells.Item(27, 8).Value = 98.032600000000002
$ws.Cells.Item(27, 9).Value = 'CPU bound {\tt [$n$: 1] [seed: 2] [$\lambda$: 0.01] [limit: 256] [$t_{cs}$: 4] [$\alpha$: 0.5] [$t_{slice}$: 128]}'

$ws.Cells.Item(28, 1).Value = 6
$ws.Cells.Item(28, 2).Value = "SJF"
$ws.Cells.Item(28, 3).Value = 1629.22
$ws.Cells.Item(28, 4).Value = 2596.35
$ws.Cells.Item(28, 5).Value = 4229.57
$ws.Cells.Item(28, 6).Value = 23
$ws.Cells.Item(28, 7).Value = 0
$ws.Cells.Item(28, 8).Value = 94.760300000000001

$ws.Cells.Item(29, 1).Value = 6
$ws.Cells.Item(29, 2).Value = "SRT"
$ws.Cells.Item(29, 3).Value = 1629.22
$ws.Cells.Item(29, 4).Value = 2543.5700000000002
$ws.Cells.Item(29, 5).Value = 4178.17
$ws.Cells.Item(29, 6).Value = 31
$ws.Cells.Item(29, 7).Value = 8
$ws.Cells.Item(29, 8).Value = 94.683599999999998

$ws.Cells.Item(30, 1).Value = 6
$ws.Cells.Item(30, 2).Value = "RR"
$ws.Cells.Item(30, 3).Value = 1629.22
$ws.Cells.Item(30, 4).Value = 2249.39
$ws.Cells.Item(30, 5).Value = 3913.57
$ws.Cells.Item(30, 6).Value = 201
$ws.Cells.Item(30, 7).Value = 178
$ws.Cells.Item(30, 8).Value = 92.023600000000002

$ws.Cells.Item(32, 1).Value = 7
$ws.Cells.Item(32, 2).Value = "FCFS"
$ws.Cells.Item(32, 3).Value = 1629.22
$ws.Cells.Item(32, 4).Value = 3216.87
$ws.Cells.Item(32, 5).Value = 4850.09
$ws.Cells.Item(32, 6).Value = 23
$ws.Cells.Item(32, 7).Value = 0
$ws.Cells.Item(32, 8).Value = 98.032600000000002
$ws.Cells.Item(32, 9).Value = 'CPU bound {\tt [$n$: 16] [seed: 2] [$\lambda$: 0.01] [limit: 256] [$t_{cs}$: 4] [$\alpha$: 0.75] [$t_{slice}$: 64]}'

$ws.Cells.Item(33, 1).Value = 7
$ws.Cells.Item(33, 2).Value = "SJF"
$ws.Cells.Item(33, 3).Value = 1629.22
$ws.Cells.Item(33, 4).Value = 2596.35
$ws.Cells.Item(33, 5).Value = 4229.57
$ws.Cells.Item(33, 6).Value = 23
$ws.Cells.Item(33, 7).Value = 0
$ws.Cells.Item(33, 8).Value = 94.760300000000001

$ws.Cells.Item(34, 1).Value = 7
$ws.Cells.Item(34, 2).Value = "SRT"
$ws.Cells.Item(34, 3).Value = 1629.22
$ws.Cells.Item(34, 4).Value = 2543.5700000000002
$ws.Cells.Item(34, 5).Value = 4178.17
$ws.Cells.Item(34, 6).Value = 31
$ws.Cells.Item(34, 7).Value = 8
$ws.Cells.Item(34, 8).Value = 94.683599999999998

$ws.Cells.Item(35, 1).Value = 7
$ws.Cells.Item(35, 2).Value = "RR"
$ws.Cells.Item(35, 3).Value = 1629.22
$ws.Cells.Item(35, 4).Value = 2327.13
$ws.Cells.Item(35, 5).Value = 4024.17
$ws.Cells.Item(35, 6).Value = 390
$ws.Cells.Item(35, 7).Value = 367
$ws.Cells.Item(35, 8).Value = 90.068299999999994

$ws.Cells.Item(37, 1).Value = 8
$ws.Cells.Item(37, 2).Value = "FCFS"
$ws.Cells.Item(37, 3).Value = 1629.22
$ws.Cells.Item(37, 4).Value = 3216.87
$ws.Cells.Item(37, 5).Value = 4850.09
$ws.Cells.Item(37, 6).Value = 23
$ws.Cells.Item(37, 7).Value = 0
$ws.Cells.Item(37, 8).Value = 98.032600000000002
$ws.Cells.Item(37, 9).Value = 'CPU bound {\tt [$n$: 8] [seed: 64] [$\lambda$: 0.001] [limit: 4096] [$t_{cs}$: 4] [$\alpha$: 0.5] [$t_{slice}$: 2048]}'

$ws.Cells.Item(38, 1).Value = 8
$ws.Cells.Item(38, 2).Value = "SJF"
$ws.Cells.Item(38, 3).Value = 1629.22
$ws.Cells.Item(38, 4).Value = 2730.65
$ws.Cells.Item(38, 5).Value = 4363.87
$ws.Cells.Item(38, 6).Value = 23
$ws.Cells.Item(38, 7).Value = 0
$ws.Cells.Item(38, 8).Value = 95.314599999999999

$ws.Cells.Item(39, 1).Value = 8
$ws.Cells.Item(39, 2).Value = "SRT"
$ws.Cells.Item(39, 3).Value = 1629.22
$ws.Cells.Item(39, 4).Value = 2017.39
$ws.Cells.Item(39, 5).Value = 3652.35
$ws.Cells.Item(39, 6).Value = 33
$ws.Cells.Item(39, 7).Value = 10
$ws.Cells.Item(39, 8).Value = 94.664500000000004

$ws.Cells.Item(40, 1).Value = 8
$ws.Cells.Item(40, 2).Value = "RR"
$ws.Cells.Item(40, 3).Value = 1629.22
$ws.Cells.Item(40, 4).Value = 2829.35
$ws.Cells.Item(40, 5).Value = 4463.78
$ws.Cells.Item(40, 6).Value = 30
$ws.Cells.Item(40, 7).Value = 7
$ws.Cells.Item(40, 8).Value = 95.246799999999993

$ws.Cells.Item(42, 1).Value = 9
$ws.Cells.Item(42, 2).Value = "FCFS"
$ws.Cells.Item(42, 3).Value = 1629.22
$ws.Cells.Item(42, 4).Value = 3216.87
$ws.Cells.Item(42, 5).Value = 4850.09
$ws.Cells.Item(42, 6).Value = 23
$ws.Cells.Item(42, 7).Value = 0
$ws.Cells.Item(42, 8).Value = 98.032600000000002
$ws.Cells.Item(42, 9).Value = 'CPU bound {\tt [$n$: 8] [seed: 64] [$\lambda$: 0.001] [limit: 4096] [$t_{cs}$: 4] [$\alpha$: 0.5] [$t_{slice}$: 2048][$rr_{add}$: BEGINNING]}'

$ws.Cells.Item(43, 1).Value = 9
$ws.Cells.Item(43, 2).Value = "SJF"
$ws.Cells.Item(43, 3).Value = 1629.22
$ws.Cells.Item(43, 4).Value = 2730.65
$ws.Cells.Item(43, 5).Value = 4363.87
$ws.Cells.Item(43, 6).Value = 23
$ws.Cells.Item(43, 7).Value = 0
$ws.Cells.Item(43, 8).Value = 95.314599999999999

$ws.Cells.Item(44, 1).Value = 9
$ws.Cells.Item(44, 2).Value = "SRT"
$ws.Cells.Item(44, 3).Value = 1629.22
$ws.Cells.Item(44, 4).Value = 2017.39
$ws.Cells.Item(44, 5).Value = 3652.35
$ws.Cells.Item(44, 6).Value = 33
$ws.Cells.Item(44, 7).Value = 10
$ws.Cells.Item(44, 8).Value = 94.664500000000004

$ws.Cells.Item(45, 1).Value = 9
$ws.Cells.Item(45, 2).Value = "RR"
$ws.Cells.Item(45, 3).Value = 1629.22
$ws.Cells.Item(45, 4).Value = 1743.61
$ws.Cells.Item(45, 5).Value = 3378.04
$ws.Cells.Item(45, 6).Value = 30
$ws.Cells.Item(45, 7).Value = 7
$ws.Cells.Item(45, 8).Value = 95.616200000000006

$ws.Cells.Item(47, 1).Value = 10
$ws.Cells.Item(47, 2).Value = "FCFS"
$ws.Cells.Item(47, 3).Value = 1629.22
$ws.Cells.Item(47, 4).Value = 3255.13
$ws.Cells.Item(47, 5).Value = 4904.3500000000004
$ws.Cells.Item(47, 6).Value = 23
$ws.Cells.Item(47, 7).Value = 0
$ws.Cells.Item(47, 8).Value = 97.097800000000007
$ws.Cells.Item(47, 9).Value = 'CPU bound {\tt [$n$: 8] [seed: 64] [$\lambda$: 0.001] [limit: 4096] [$t_{cs}$: 20] [$\alpha$: 0.5] [$t_{slice}$: 2048]}'

$ws.Cells.Item(48, 1).Value = 10
$ws.Cells.Item(48, 2).Value = "SJF"
$ws.Cells.Item(48, 3).Value = 1629.22
$ws.Cells.Item(48, 4).Value = 2761.26
$ws.Cells.Item(48, 5).Value = 4410.4799999999996
$ws.Cells.Item(48, 6).Value = 23
$ws.Cells.Item(48, 7).Value = 0
$ws.Cells.Item(48, 8).Value = 94.430700000000002

$ws.Cells.Item(49, 1).Value = 10
$ws.Cells.Item(49, 2).Value = "SRT"
$ws.Cells.Item(49, 3).Value = 1629.22
$ws.Cells.Item(49, 4).Value = 2057.7399999999998
$ws.Cells.Item(49, 5).Value = 3715.65
$ws.Cells.Item(49, 6).Value = 33
$ws.Cells.Item(49, 7).Value = 10
$ws.Cells.Item(49, 8).Value = 93.418400000000005

$ws.Cells.Item(50, 1).Value = 10
$ws.Cells.Item(50, 2).Value = "RR"
$ws.Cells.Item(50, 3).Value = 1629.22
$ws.Cells.Item(50, 4).Value = 2876.65
$ws.Cells.Item(50, 5).Value = 4531.96
$ws.Cells.Item(50, 6).Value = 30
$ws.Cells.Item(50, 7).Value = 7
$ws.Cells.Item(50, 8).Value = 94.098699999999994

# Update selection to match the final active cell
$ws.Range("I47").Select()
